$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Nid2"
$ws.Range("C2").Value = "Col13a1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 18.904217
$ws.Range("H2").Value = 56.712651
$ws.Range("I2").Value = 0.3376032603366536
$ws.Range("J2").Value = 0.3376032603366536
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.3943663333333334
$ws.Range("N2").Value = 1.183099
$ws.Range("O2").Value = 0.393328993667368
$ws.Range("P2").Value = 0.393328993667368
$ws.Range("Q2").Value = 7.455186742827667
$ws.Range("R2").Value = 67.096680685449
$ws.Range("S2").Value = 0.1327891506470384
$ws.Range("T2").Value = 0.1327891506470384

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Nid2"
$ws.Range("C3").Value = "Col13a1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 18.904217
$ws.Range("H3").Value = 56.712651
$ws.Range("I3").Value = 0.3376032603366536
$ws.Range("J3").Value = 0.3376032603366536
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.043492
$ws.Range("N3").Value = 0.130476
$ws.Range("O3").Value = 0.04337759881273123
$ws.Range("P3").Value = 0.04337759881273123
$ws.Range("Q3").Value = 0.822182205764
$ws.Range("R3").Value = 7.399639851876
$ws.Range("S3").Value = 0.01464441878475342
$ws.Range("T3").Value = 0.01464441878475342

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Nid2"
$ws.Range("C4").Value = "Col13a1"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 18.904217
$ws.Range("H4").Value = 56.712651
$ws.Range("I4").Value = 0.3376032603366536
$ws.Range("J4").Value = 0.3376032603366536
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.564779
$ws.Range("N4").Value = 1.694337
$ws.Range("O4").Value = 0.5632934075199009
$ws.Range("P4").Value = 0.5632934075199009
$ws.Range("Q4").Value = 10.676704773043
$ws.Range("R4").Value = 96.090342957387
$ws.Range("S4").Value = 0.1901696909048618
$ws.Range("T4").Value = 0.1901696909048618

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Nid2"
$ws.Range("C5").Value = "Col13a1"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 21.42159766666667
$ws.Range("H5").Value = 64.264793
$ws.Range("I5").Value = 0.3825602093906729
$ws.Range("J5").Value = 0.3825602093906729
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.3943663333333334
$ws.Range("N5").Value = 1.183099
$ws.Range("O5").Value = 0.393328993667368
$ws.Range("P5").Value = 0.393328993667368
$ws.Range("Q5").Value = 8.447956925945224
$ws.Range("R5").Value = 76.031612333507
$ws.Range("S5").Value = 0.1504720221768109
$ws.Range("T5").Value = 0.1504720221768109

$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Nid2"
$ws.Range("C6").Value = "Col13a1"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 21.42159766666667
$ws.Range("H6").Value = 64.264793
$ws.Range("I6").Value = 0.3825602093906729
$ws.Range("J6").Value = 0.3825602093906729
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.043492
$ws.Range("N6").Value = 0.130476
$ws.Range("O6").Value = 0.04337759881273123
$ws.Range("P6").Value = 0.04337759881273123
$ws.Range("Q6").Value = 0.9316681257186668
$ws.Range("R6").Value = 8.385013131468
$ws.Range("S6").Value = 0.01659454328466306
$ws.Range("T6").Value = 0.01659454328466306

$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Nid2"
$ws.Range("C7").Value = "Col13a1"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 21.42159766666667
$ws.Range("H7").Value = 64.264793
$ws.Range("I7").Value = 0.3825602093906729
$ws.Range("J7").Value = 0.3825602093906729
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.564779
$ws.Range("N7").Value = 1.694337
$ws.Range("O7").Value = 0.5632934075199009
$ws.Range("P7").Value = 0.5632934075199009
$ws.Range("Q7").Value = 12.09846850858233
$ws.Range("R7").Value = 108.886216577241
$ws.Range("S7").Value = 0.2154936439291989
$ws.Range("T7").Value = 0.2154936439291989

$ws.Range("A8").Value = "M1"
$ws.Range("B8").Value = "Nid2"
$ws.Range("C8").Value = "Col13a1"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.6843913333333332
$ws.Range("H8").Value = 2.053174
$ws.Range("I8").Value = 0.01222228593120163
$ws.Range("J8").Value = 0.01222228593120164
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.3943663333333334
$ws.Range("N8").Value = 1.183099
$ws.Range("O8").Value = 0.393328993667368
$ws.Range("P8").Value = 0.393328993667368
$ws.Range("Q8").Value = 0.2699009006917777
$ws.Range("R8").Value = 2.429108106226
$ws.Range("S8").Value = 0.004807379425634368
$ws.Range("T8").Value = 0.004807379425634369

$ws.Range("A9").Value = "M1"
$ws.Range("B9").Value = "Nid2"
$ws.Range("C9").Value = "Col13a1"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.6843913333333332
$ws.Range("H9").Value = 2.053174
$ws.Range("I9").Value = 0.01222228593120163
$ws.Range("J9").Value = 0.01222228593120164
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.043492
$ws.Range("N9").Value = 0.130476
$ws.Range("O9").Value = 0.04337759881273123
$ws.Range("P9").Value = 0.04337759881273123
$ws.Range("Q9").Value = 0.02976554786933333
$ws.Range("R9").Value = 0.267889930824
$ws.Range("S9").Value = 0.0005301734156981536
$ws.Range("T9").Value = 0.0005301734156981536

$ws.Range("A10").Value = "M1"
$ws.Range("B10").Value = "Nid2"
$ws.Range("C10").Value = "Col13a1"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.6843913333333332
$ws.Range("H10").Value = 2.053174
$ws.Range("I10").Value = 0.01222228593120163
$ws.Range("J10").Value = 0.01222228593120164
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.564779
$ws.Range("N10").Value = 1.694337
$ws.Range("O10").Value = 0.5632934075199009
$ws.Range("P10").Value = 0.5632934075199009
$ws.Range("Q10").Value = 0.3865298528486666
$ws.Range("R10").Value = 3.478768675638
$ws.Range("S10").Value = 0.006884733089869114
$ws.Range("T10").Value = 0.006884733089869114

$ws.Range("A11").Value = "M2"
$ws.Range("B11").Value = "Nid2"
$ws.Range("C11").Value = "Col13a1"
$ws.Range("D11").Value = "ECs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.4048903333333334
$ws.Range("H11").Value = 1.214671
$ws.Range("I11").Value = 0.007230783301531494
$ws.Range("J11").Value = 0.007230783301531494
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 0.3943663333333334
$ws.Range("N11").Value = 1.183099
$ws.Range("O11").Value = 0.393328993667368
$ws.Range("P11").Value = 0.393328993667368
$ws.Range("Q11").Value = 0.1596751161587778
$ws.Range("R11").Value = 1.437076045429
$ws.Range("S11").Value = 0.002844076719418191
$ws.Range("T11").Value = 0.002844076719418191

$ws.Range("A12").Value = "M2"
$ws.Range("B12").Value = "Nid2"
$ws.Range("C12").Value = "Col13a1"
$ws.Range("D12").Value = "FAPs"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 0.4048903333333334
$ws.Range("H12").Value = 1.214671
$ws.Range("I12").Value = 0.007230783301531494
$ws.Range("J12").Value = 0.007230783301531494
$ws.Range("K12").Value = 1
$ws.Range("L12").Value = 0.3333333333333333
$ws.Range("M12").Value = 0.043492
$ws.Range("N12").Value = 0.130476
$ws.Range("O12").Value = 0.04337759881273123
$ws.Range("P12").Value = 0.04337759881273123
$ws.Range("Q12").Value = 0.01760949037733334
$ws.Range("R12").Value = 0.158485413396
$ws.Range("S12").Value = 0.0003136540171556293
$ws.Range("T12").Value = 0.0003136540171556293

$ws.Range("A13").Value = "M2"
$ws.Range("B13").Value = "Nid2"
$ws.Range("C13").Value = "Col13a1"
$ws.Range("D13").Value = "sCs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 0.4048903333333334
$ws.Range("H13").Value = 1.214671
$ws.Range("I13").Value = 0.007230783301531494
$ws.Range("J13").Value = 0.007230783301531494
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.564779
$ws.Range("N13").Value = 1.694337
$ws.Range("O13").Value = 0.5632934075199009
$ws.Range("P13").Value = 0.5632934075199009
$ws.Range("Q13").Value = 0.2286735575696667
$ws.Range("R13").Value = 2.058062018127
$ws.Range("S13").Value = 0.004073052564957674
$ws.Range("T13").Value = 0.004073052564957674

$ws.Range("A14").Value = "Neutro"
$ws.Range("B14").Value = "Nid2"
$ws.Range("C14").Value = "Col13a1"
$ws.Range("D14").Value = "ECs"
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 0.525581
$ws.Range("H14").Value = 1.576743
$ws.Range("I14").Value = 0.009386152262799286
$ws.Range("J14").Value = 0.009386152262799286
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 0.3943663333333334
$ws.Range("N14").Value = 1.183099
$ws.Range("O14").Value = 0.393328993667368
$ws.Range("P14").Value = 0.393328993667368
$ws.Range("Q14").Value = 0.2072714518396667
$ws.Range("R14").Value = 1.865443066557
$ws.Range("S14").Value = 0.003691845823935532
$ws.Range("T14").Value = 0.003691845823935532

$ws.Range("A15").Value = "Neutro"
$ws.Range("B15").Value = "Nid2"
$ws.Range("C15").Value = "Col13a1"
$ws.Range("D15").Value = "FAPs"
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 0.525581
$ws.Range("H15").Value = 1.576743
$ws.Range("I15").Value = 0.009386152262799286
$ws.Range("J15").Value = 0.009386152262799286
$ws.Range("K15").Value = 1
$ws.Range("L15").Value = 0.3333333333333333
$ws.Range("M15").Value = 0.043492
$ws.Range("N15").Value = 0.130476
$ws.Range("O15").Value = 0.04337759881273123
$ws.Range("P15").Value = 0.04337759881273123
$ws.Range("Q15").Value = 0.022858568852
$ws.Range("R15").Value = 0.205727119668
$ws.Range("S15").Value = 0.0004071487472509168
$ws.Range("T15").Value = 0.0004071487472509168

$ws.Range("A16").Value = "Neutro"
$ws.Range("B16").Value = "Nid2"
$ws.Range("C16").Value = "Col13a1"
$ws.Range("D16").Value = "sCs"
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 0.525581
$ws.Range("H16").Value = 1.576743
$ws.Range("I16").Value = 0.009386152262799286
$ws.Range("J16").Value = 0.009386152262799286
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 0.564779
$ws.Range("N16").Value = 1.694337
$ws.Range("O16").Value = 0.5632934075199009
$ws.Range("P16").Value = 0.5632934075199009
$ws.Range("Q16").Value = 0.296837111599
$ws.Range("R16").Value = 2.671534004391
$ws.Range("S16").Value = 0.005287157691612838
$ws.Range("T16").Value = 0.005287157691612838

$ws.Range("A17").Value = "sCs"
$ws.Range("B17").Value = "Nid2"
$ws.Range("C17").Value = "Col13a1"
$ws.Range("D17").Value = "ECs"
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 14.05468533333333
$ws.Range("H17").Value = 42.164056
$ws.Range("I17").Value = 0.2509973087771411
$ws.Range("J17").Value = 0.2509973087771412
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 0.3943663333333334
$ws.Range("N17").Value = 1.183099
$ws.Range("O17").Value = 0.393328993667368
$ws.Range("P17").Value = 0.393328993667368
$ws.Range("Q17").Value = 5.542694721060445
$ws.Range("R17").Value = 49.88425248954401
$ws.Range("S17").Value = 0.09872451887453053
$ws.Range("T17").Value = 0.09872451887453056

$ws.Range("A18").Value = "sCs"
$ws.Range("B18").Value = "Nid2"
$ws.Range("C18").Value = "Col13a1"
$ws.Range("D18").Value = "FAPs"
$ws.Range("E18").Value = 3
$ws.Range("F18").Value = 1
$ws.Range("G18").Value = 14.05468533333333
$ws.Range("H18").Value = 42.164056
$ws.Range("I18").Value = 0.2509973087771411
$ws.Range("J18").Value = 0.2509973087771412
$ws.Range("K18").Value = 1
$ws.Range("L18").Value = 0.3333333333333333
$ws.Range("M18").Value = 0.043492
$ws.Range("N18").Value = 0.130476
$ws.Range("O18").Value = 0.04337759881273123
$ws.Range("P18").Value = 0.04337759881273123
$ws.Range("Q18").Value = 0.6112663745173333
$ws.Range("R18").Value = 5.501397370656001
$ws.Range("S18").Value = 0.01088766056321005
$ws.Range("T18").Value = 0.01088766056321005

$ws.Range("A19").Value = "sCs"
$ws.Range("B19").Value = "Nid2"
$ws.Range("C19").Value = "Col13a1"
$ws.Range("D19").Value = "sCs"
$ws.Range("E19").Value = 3
$ws.Range("F19").Value = 1
$ws.Range("G19").Value = 14.05468533333333
$ws.Range("H19").Value = 42.164056
$ws.Range("I19").Value = 0.2509973087771411
$ws.Range("J19").Value = 0.2509973087771412
$ws.Range("K19").Value = 3
$ws.Range("L19").Value = 1
$ws.Range("M19").Value = 0.564779
$ws.Range("N19").Value = 1.694337
$ws.Range("O19").Value = 0.5632934075199009
$ws.Range("P19").Value = 0.5632934075199009
$ws.Range("Q19").Value = 7.937791127874667
$ws.Range("R19").Value = 71.440120150872
$ws.Range("S19").Value = 0.1413851293394005
$ws.Range("T19").Value = 0.1413851293394006

